$p = $ppt.ActivePresentation

# Add the new "Ejercicio práctico" slide at the end of the deck, using the
# same "Title and Content" layout (index 2) as the rest of the deck's slides.
$slideIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($slideIndex, 2)

# --- Title placeholder -------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Name = "Título 1"
$title.TextFrame.TextRange.Text = "Ejercicio práctico "
$title.TextFrame.TextRange.LanguageID = "es-ES"

# --- Body / content placeholder -----------------------------------------
$body = $s.Shapes.Item(2)
$body.Name = "Marcador de contenido 2"

# Explicit size/position override, matching the authored slide.
$body.Left = 108.0
$body.Top = 130.1284451968504
$body.Width = 756.0
$body.Height = 391.04594551181106

# "Shrink text on overflow" autofit (renders as <a:normAutofit/>).
$body.TextFrame.AutoSize = 2

$tf = $body.TextFrame

# Paragraph 1
$tf.TextRange.Text = "Crear una nueva aplicación API-REST utilizando las plantillas de visual studio llamada: Pract1"
$tf.TextRange.LanguageID = "es-ES"
$p1 = $tf.TextRange.Paragraphs(1, 1)
$p1.Characters($p1.Text.Length - "Pract1".Length + 1, "Pract1".Length).Font.Italic = $true

# Paragraph 2
[void]$tf.TextRange.InsertAfter("`rAñadir un nuevo fichero de configuración .json en el siguiente directorio: Infrastructure/pract1SettingDemo.json")
$n = $tf.TextRange.Paragraphs().Count
$p2 = $tf.TextRange.Paragraphs($n, 1)
$p2.LanguageID = "es-ES"
$italic2 = "Infrastructure/pract1SettingDemo.json"
$p2.Characters($p2.Text.Length - $italic2.Length + 1, $italic2.Length).Font.Italic = $true

# Paragraph 3
[void]$tf.TextRange.InsertAfter("`rRellenar el fichero con un objeto json válido ")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraph 4
[void]$tf.TextRange.InsertAfter("`rAñadir ese json al objeto IConfiguration de .NET Core con la extensión AddJson ya comentada.")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraph 5
[void]$tf.TextRange.InsertAfter("`rCrear una sección de un fichero json que puedas mapear a un objeto C# tal y como hemos visto en la demo.")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraph 6
[void]$tf.TextRange.InsertAfter("`rCrear un controlador nuevo en esa API-REST con el nombre: TestController")
$n = $tf.TextRange.Paragraphs().Count
$p6 = $tf.TextRange.Paragraphs($n, 1)
$p6.LanguageID = "es-ES"
$italic6 = "TestController"
$p6.Characters($p6.Text.Length - $italic6.Length + 1, $italic6.Length).Font.Italic = $true

# Paragraph 7
[void]$tf.TextRange.InsertAfter("`rAñadir una ruta nueva en ese controlador que empiece por: /api/nombrecontrolador")
$n = $tf.TextRange.Paragraphs().Count
$p7 = $tf.TextRange.Paragraphs($n, 1)
$p7.LanguageID = "es-ES"
$italic7 = "/api/nombrecontrolador"
$p7.Characters($p7.Text.Length - $italic7.Length + 1, $italic7.Length).Font.Italic = $true

# Paragraph 8
[void]$tf.TextRange.InsertAfter("`rCrear una variable de entorno nueva dentro de visual studio.")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraph 9
[void]$tf.TextRange.InsertAfter("`rCrear una acción del controlador que devuelva el valor de una sección del fichero json mediante un objeto C#")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraph 10
[void]$tf.TextRange.InsertAfter("`rCrear una acción del controlador que devuelva el valor de una propiedad del json mediante el objeto IConfiguration")
$n = $tf.TextRange.Paragraphs().Count
$p10 = $tf.TextRange.Paragraphs($n, 1)
$p10.LanguageID = "es-ES"
$italic10 = "IConfiguration"
$p10.Characters($p10.Text.Length - $italic10.Length + 1, $italic10.Length).Font.Italic = $true

# Paragraph 11
[void]$tf.TextRange.InsertAfter("`rCrear una acción del controlador que devuelva el valor de la variable de entorno anteriormente comentada.")
$n = $tf.TextRange.Paragraphs().Count
$tf.TextRange.Paragraphs($n, 1).LanguageID = "es-ES"

# Paragraphs 12 and 13 (trailing blank lines)
[void]$tf.TextRange.InsertAfter("`r")
[void]$tf.TextRange.InsertAfter("`r")
